$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the date strings in column A (slashes -> dashes) -----------
# Plain assignment works for most rows, but some day/month combos are
# ambiguous (e.g. "01-08-2022" could be parsed as a US-style M-D-Y date)
# and Excel's COM layer would silently convert them to a date serial
# number instead of keeping them as text. For those rows we briefly force
# a Text number format, assign the value, then restore the style to
# "Normal" so the cell ends up with plain text and no residual
# style/format change.

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("A3")  "28-07-2022"
Set-TextValue $ws.Range("A4")  "01-08-2022"
Set-TextValue $ws.Range("A5")  "04-08-2022"
Set-TextValue $ws.Range("A6")  "08-08-2022"
Set-TextValue $ws.Range("A7")  "11-08-2022"
Set-TextValue $ws.Range("A8")  "15-08-2022"
Set-TextValue $ws.Range("A9")  "18-08-2022"
Set-TextValue $ws.Range("A10") "22-08-2022"
Set-TextValue $ws.Range("A11") "25-08-2022"
Set-TextValue $ws.Range("A12") "29-08-2022"
Set-TextValue $ws.Range("A13") "01-09-2022"
Set-TextValue $ws.Range("A14") "05-09-2022"
Set-TextValue $ws.Range("A15") "08-09-2022"
Set-TextValue $ws.Range("A16") "12-09-2022"
Set-TextValue $ws.Range("A17") "15-09-2022"
Set-TextValue $ws.Range("A18") "19-09-2022"
Set-TextValue $ws.Range("A19") "22-09-2022"
Set-TextValue $ws.Range("A20") "26-09-2022"
Set-TextValue $ws.Range("A21") "29-09-2022"

# --- Update the attendance counts that changed --------------------------

# Row 3: D 0->1, G 0->1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: D 0->2, E 0->1, F 0->1, H 1->0
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = 0

# Row 6: D 0->1, E 0->1, H 1->0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 13: D 0->1, E 0->1, H 1->0
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14: D 0->1, E 0->1, H 1->0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0
